$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title paragraph: the trailing run of spaces grows from 130 to 131
#    characters because the single-space run that used to sit in front of
#    the (hidden) "_GoBack" bookmark is folded into it.
#    ("Input Format Guidelines" + " " + 130 spaces  ->  "Input Format
#    Guidelines" + 131 spaces)
# ---------------------------------------------------------------------------
$titleRng = $d.Content
$titleHit = $titleRng.Find.Execute("Input Format Guidelines", $true, $false,
                                    $false, $false, $false, $true, 1, $false,
                                    "", 0)
if ($titleHit) {
    $afterTitle = $titleRng.End

    # locate the lone space run right after the title text
    $spaceRun = $d.Range($afterTitle, $afterTitle + 1)
    if ($spaceRun.Text -eq " ") {
        $spaceRun.Text = ""
    }

    # the 130-space run now starts where the single space run used to be;
    # turn its first character into two characters so the run grows by one
    # space (130 -> 131 total)
    $firstOfRun = $d.Range($afterTitle, $afterTitle + 1)
    if ($firstOfRun.Text -eq " ") {
        $firstOfRun.Text = "  "
    }
}

# ---------------------------------------------------------------------------
# 2) Fix typo "can put entered" -> "can be entered"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("can put entered", $true, $false, $false, $false,
                         $false, $true, 1, $false, "can be entered", 2)

# ---------------------------------------------------------------------------
# 3) Reword the multiplication-operator bullet
# ---------------------------------------------------------------------------
$oldMult = "Multiplication operators must be used wherever they are implied. "
$newMult = "The " + [char]0x201C + "*" + [char]0x201D +
           " operator must be used wherever multiplication is implied. "
$d.Content.Find.Execute($oldMult, $true, $false, $false, $false, $false,
                         $true, 1, $false, $newMult, 2)

Write-Output "done"
